$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arm")

# Clear the raw data in J4 (Init_C for "mg") as part of the data refresh
$ws.Range("J4").ClearContents()

# Leave the active selection on the cleared cell, matching the edit session
$ws.Range("J4").Select()
